$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("FEINmismatch")
$sheet2 = $wb.Worksheets.Item("FEINSSNmismatch")

$sheet1Dates = @(
    "Sat Oct 12 23:20:28 EDT 2024",
    "Sat Oct 12 23:20:43 EDT 2024",
    "Sat Oct 12 23:20:55 EDT 2024",
    "Sat Oct 12 23:21:07 EDT 2024",
    "Sat Oct 12 23:21:21 EDT 2024",
    "Sat Oct 12 23:21:35 EDT 2024",
    "Sat Oct 12 23:21:47 EDT 2024",
    "Sat Oct 12 23:21:59 EDT 2024",
    "Sat Oct 12 23:22:10 EDT 2024",
    "Sat Oct 12 23:22:22 EDT 2024",
    "Sat Oct 12 23:22:33 EDT 2024",
    "Sat Oct 12 23:22:45 EDT 2024",
    "Sat Oct 12 23:22:57 EDT 2024",
    "Sat Oct 12 23:23:09 EDT 2024",
    "Sat Oct 12 23:23:21 EDT 2024",
    "Sat Oct 12 23:23:32 EDT 2024",
    "Sat Oct 12 23:23:44 EDT 2024",
    "Sat Oct 12 23:23:55 EDT 2024",
    "Sat Oct 12 23:24:07 EDT 2024",
    "Sat Oct 12 23:24:18 EDT 2024",
    "Sat Oct 12 23:24:30 EDT 2024",
    "Sat Oct 12 23:24:41 EDT 2024",
    "Sat Oct 12 23:24:53 EDT 2024",
    "Sat Oct 12 23:25:04 EDT 2024",
    "Sat Oct 12 23:25:16 EDT 2024",
    "Sat Oct 12 23:25:28 EDT 2024",
    "Sat Oct 12 23:25:39 EDT 2024",
    "Sat Oct 12 23:25:51 EDT 2024",
    "Sat Oct 12 23:26:04 EDT 2024"
)

$sheet2Dates = @(
    "Sat Oct 12 23:26:16 EDT 2024",
    "Sat Oct 12 23:26:28 EDT 2024",
    "Sat Oct 12 23:26:39 EDT 2024",
    "Sat Oct 12 23:26:50 EDT 2024",
    "Sat Oct 12 23:27:02 EDT 2024",
    "Sat Oct 12 23:27:13 EDT 2024",
    "Sat Oct 12 23:27:24 EDT 2024",
    "Sat Oct 12 23:27:35 EDT 2024",
    "Sat Oct 12 23:27:47 EDT 2024",
    "Sat Oct 12 23:27:58 EDT 2024",
    "Sat Oct 12 23:28:10 EDT 2024",
    "Sat Oct 12 23:28:22 EDT 2024",
    "Sat Oct 12 23:28:33 EDT 2024",
    "Sat Oct 12 23:28:44 EDT 2024",
    "Sat Oct 12 23:28:56 EDT 2024",
    "Sat Oct 12 23:29:08 EDT 2024",
    "Sat Oct 12 23:29:19 EDT 2024",
    "Sat Oct 12 23:29:30 EDT 2024"
)

for ($i = 0; $i -lt $sheet1Dates.Length; $i++) {
    $row = $i + 2
    $sheet1.Cells.Item($row, 2).Value = $sheet1Dates[$i]
}

for ($i = 0; $i -lt $sheet2Dates.Length; $i++) {
    $row = $i + 2
    $sheet2.Cells.Item($row, 2).Value = $sheet2Dates[$i]
}
